$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old numeric data rows (A1:A3 held 1137/1603/1907)
$ws.Range("A1:A3").ClearContents()

# Write the new header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Status"
$ws.Range("C1").Value = "Date&Time"
